# repull data, push all data, mean calculation
# Update column F ("dSF") values for rows 3-35 (skipping rows already correct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    4  = 1
    5  = 2
    6  = 6
    7  = 7
    8  = -3
    9  = 5
    10 = -2
    11 = 2
    13 = -3
    14 = 1
    15 = 2
    16 = 3
    18 = -2
    19 = -1
    20 = 4
    21 = 7
    22 = 3
    23 = 5
    24 = 1
    25 = 4
    26 = 3
    27 = 2
    28 = 2
    29 = 6
    30 = -5
    31 = -1
    32 = 1
    33 = -4
    34 = 1
    35 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
